# Move the last data row (row 5: ID 444 / "Guest") down to row 8,
# leaving rows 6-7 blank, and update the selection to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5:B5").Cut($ws.Range("A8")) | Out-Null
$ws.Range("A8:B8").Select() | Out-Null
